$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to remain TEXT (avoids Excel
# auto-converting number-looking strings to numeric cells / losing precision
# or trailing zeros), and without leaving a permanent style/format change on
# the cell (ClearFormats resets format back to the workbook default after).
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = '36.565.75'
$ws.Range("E2").Value = '  -2.49%  '

# Row 3
$ws.Range("D3").Value = '1.985.71'
$ws.Range("E3").Value = '  -3.51%  '

# Row 4
$ws.Range("E4").Value = '  +0.17%  '

# Row 5
Set-TextValue $ws.Range("D5") '245.25'
$ws.Range("E5").Value = '  -2.81%  '

# Row 6
Set-TextValue $ws.Range("D6") '0.628'
$ws.Range("E6").Value = '  -3.33%  '

# Row 7
Set-TextValue $ws.Range("D7") '58.77'
$ws.Range("E7").Value = '  -11.37%  '

# Row 8
$ws.Range("E8").Value = '  +0.04%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.377'
$ws.Range("E9").Value = '  -1.49%  '

# Row 10
Set-TextValue $ws.Range("D10") '57.68'
$ws.Range("E10").Value = '  -3.25%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0816'
$ws.Range("E11").Value = '  +6.09%  '

# Row 13
Set-TextValue $ws.Range("D13") '24.38'
$ws.Range("E13").Value = '  +14.15%  '

# Row 14
Set-TextValue $ws.Range("D14") '0.869'
$ws.Range("E14").Value = '  -4.61%  '

# Row 15
Set-TextValue $ws.Range("D15") '14.18'
$ws.Range("E15").Value = '  -4.91%  '

# Row 16
$ws.Range("D16").Value = '2.276.96'
$ws.Range("E16").Value = '  -3.52%  '

# Row 17
Set-TextValue $ws.Range("D17") '5.46'
$ws.Range("E17").Value = '  -2.27%  '

# Row 18
$ws.Range("D18").Value = '1.985.54'
$ws.Range("E18").Value = '  -3.34%  '

# Row 19
$ws.Range("D19").Value = '36.459.56'
$ws.Range("E19").Value = '  -2.08%  '

# Row 20
Set-TextValue $ws.Range("D20") '71.21'
$ws.Range("E20").Value = '  -3.73%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0866'
$ws.Range("E21").Value = '  -1.53%  '

# Row 22
Set-TextValue $ws.Range("D22") '5.34'
$ws.Range("E22").Value = '  -2.65%  '

# Row 23
Set-TextValue $ws.Range("D23") '235.39'
$ws.Range("E23").Value = '  -1.96%  '

# Row 24
Set-TextValue $ws.Range("D24") '2.63'
$ws.Range("E24").Value = '  -1.16%  '

# Row 25
$ws.Range("E25").Value = '  +0.10%  '

# Row 26
Set-TextValue $ws.Range("D26") '2.31'
$ws.Range("E26").Value = '  -3.11%  '

# Row 27
Set-TextValue $ws.Range("D27") '10.17'
$ws.Range("E27").Value = '  +4.60%  '

# Row 28
Set-TextValue $ws.Range("D28") '162.34'
$ws.Range("E28").Value = '  +0.41%  '

# Row 29
Set-TextValue $ws.Range("D29") '19.91'
$ws.Range("E29").Value = '  -0.58%  '

# Row 30
Set-TextValue $ws.Range("D30") '0.128'
$ws.Range("E30").Value = '  +11.83%  '

# Row 31
Set-TextValue $ws.Range("D31") '0.120'
$ws.Range("E31").Value = '  -1.52%  '

# Row 32
Set-TextValue $ws.Range("D32") '1.18'
$ws.Range("E32").Value = '  -1.70%  '

# Row 33
Set-TextValue $ws.Range("D33") '4.94'
$ws.Range("E33").Value = '  -6.52%  '

# Row 34
$ws.Range("E34").Value = '  +1.28%  '

# Row 35
Set-TextValue $ws.Range("D35") '4.48'
$ws.Range("E35").Value = '  -6.15%  '

# Row 36
Set-TextValue $ws.Range("D36") '6.20'
$ws.Range("E36").Value = '  +1.75%  '

# Row 37
Set-TextValue $ws.Range("D37") '2.30'
$ws.Range("E37").Value = '  -6.36%  '

# Row 38
$ws.Range("E38").Value = '  +0.25%  '

# Row 39
Set-TextValue $ws.Range("D39") '1.77'
$ws.Range("E39").Value = '  -4.04%  '

# Row 40
Set-TextValue $ws.Range("D40") '3.08'
$ws.Range("E40").Value = '  +0.32%  '

# Row 41
$ws.Range("E41").Value = '  +0.39%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.0972'
$ws.Range("E42").Value = '  -6.05%  '

# Row 43
$ws.Range("E43").Value = '  -3.47%  '

# Row 44
Set-TextValue $ws.Range("D44") '0.0215'
$ws.Range("E44").Value = '  -2.56%  '

# Row 45
$ws.Range("E45").Value = '  -4.60%  '

# Row 46
Set-TextValue $ws.Range("D46") '16.34'
$ws.Range("E46").Value = '  -4.09%  '

# Row 47
Set-TextValue $ws.Range("D47") '93.00'
$ws.Range("E47").Value = '  -2.75%  '

# Row 48
Set-TextValue $ws.Range("D48") '7.65'
$ws.Range("E48").Value = '  -3.70%  '

# Row 49
$ws.Range("D49").Value = '1.376.30'
$ws.Range("E49").Value = '  -2.95%  '

# Row 50
Set-TextValue $ws.Range("D50") '2.86'
$ws.Range("E50").Value = '  -2.87%  '

# Row 51
Set-TextValue $ws.Range("D51") '45.48'
$ws.Range("E51").Value = '  -2.96%  '
